$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = 1.445647641019636;  C = 0.3048912486333797; D = 3.223369029078222;  E = 13.86384647080068;  G = 18.83775438953192 }
    3 = @{ B = 0.1169995834814548; C = 0.3048912486333797; D = 0.7210945179870265; E = 0.5333859586016987;  G = 1.67637130870356  }
    4 = @{ B = 3.272327238179451;  C = 1.626987699542094;  D = 0.1496068669990043; E = 0.5333859586016987;  G = 5.582307763322248 }
    5 = @{ B = 1.445647641019636;  C = 1.626987699542094;  D = 0.7210945179870265; E = 0.5333859586016987;  G = 4.327115817150455 }
    6 = @{ B = 1.445647641019636;  C = 1.626987699542094;  D = 3.223369029078222;  E = 0.5333859586016987;  G = 6.82939032824165  }
    7 = @{ B = 3.272327238179451;  C = 0.3048912486333797; D = 0.1496068669990043; E = 0.5333859586016987;  G = 4.260211312413533 }
    8 = @{ B = 1.445647641019636;  C = 1.626987699542094;  D = 0.7210945179870265; E = 0.5333859586016987;  G = 4.327115817150455 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
